$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14:I14").Copy()
$ws.Range("A15:I15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A15").Value = 2311
$ws.Range("B15").Value = "Longest Binary Subsequence Less Than or Equal to K"
$ws.Range("C15").Value = "#string #dp #greedy #memoization"
$ws.Range("D15").Value = "medium"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 13
$ws.Range("H15").Value = 45834
$ws.Range("I15").Value = 45834

$ws.Rows.Item(15).RowHeight = 51

$ws.Range("F18").Select() | Out-Null

try { $excel.ActiveWindow.Width = 26860 } catch {}
try { $excel.ActiveWindow.Height = 14540 } catch {}
